$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders two pairs of rows of data (title, timestamp, historical
# distance, uri) while the "time bucket" column (D) and the row positions
# themselves stay put:
#   - Row 2 ("Jay-Z confirmed to headline Glastonbury 2008") and
#     Row 3 ("Jay-Z's Glastonbury Wonderwall dig at Noel Gallagher") swap.
#   - Row 5 ("'Middle-aged rockers dilute Glastonbury spirit'") and
#     Row 7 ("Glastonbury ends on-going ticketing saga") swap.
# Column D ("time bucket") is identical for every row in each swapped pair,
# so touching it is not required (but swapping it too is harmless).
#
# We use Range.Copy (cell-to-cell) routed through an out-of-the-way scratch
# cell instead of plain Value assignment, because one of the titles
# ("'Middle-aged rockers dilute Glastonbury spirit'") begins with a literal
# apostrophe; assigning such a string straight into .Value/.Value2 makes
# Excel treat the apostrophe as a "quote prefix" marker (stripping it from
# the stored text and tagging the cell with a new number format) instead of
# keeping it as literal text. Copy preserves the exact stored value/style.

function Swap-Rows($ws, $rowA, $rowB) {
    $cols = 1,2,3,4,5  # A=title, B=timestamp, C=historical distance, D=time bucket, E=uri
    foreach ($col in $cols) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $scratch = $ws.Cells.Item(1000 + $col, 1)

        $cellA.Copy($scratch)
        $cellB.Copy($cellA)
        $scratch.Copy($cellB)
        $scratch.Clear()
    }
}

Swap-Rows $ws 2 3
Swap-Rows $ws 5 7
